# Regenerate orders with updated distance/size codes.
# Pure text substitution across all cells (labels, filenames, condition
# strings, and the Distance/Size lookup columns all encode the same
# Dxx/Sxx tokens), so a global Find/Replace reproduces the diff exactly
# without disturbing cell positions, shared-string ordering elsewhere,
# or any non-text cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.Cells

# Distance codes
$cells.Replace("D64", "D69") | Out-Null
$cells.Replace("D80", "D86") | Out-Null
$cells.Replace("D51", "D55") | Out-Null

# Size code (S20 / S25 are untouched)
$cells.Replace("S30", "S31") | Out-Null
